$d = $word.ActiveDocument

# --- Step 1: fix paragraph 2 -----------------------------------------------
# Paragraph 2 currently holds the sentence split across two runs with the
# "_GoBack" bookmark sitting in between ("...null hy" | bookmark | "pothesis
# is true. "). The target collapses this into a single clean run holding the
# complete sentence, with the bookmark removed from here (it is relocated to
# a brand new trailing paragraph below). Replace just the paragraph's content
# (i.e. everything up to, but not including, the paragraph mark) so the
# existing pPr (ListParagraph / numPr) is left completely untouched.
$p2 = $d.Paragraphs.Item(2).Range
$contentRng = $d.Range($p2.Start, $p2.End - 1)
$fixedParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">P-value is the answer to the question: What is the probability that an outcome from the null distribution is bigger than what we observed when the null hypothesis is true. </w:t></w:r></w:p>'
$contentRng.InsertXML($fixedParaXml)

# --- Step 2: append the new paragraphs -------------------------------------
# New, empty spacer paragraph + bold/underlined "Central Limit Theorem (CLT)
# and t-distribution" heading + two new bulleted ("ListParagraph") items, the
# last of which now carries the relocated (empty) "_GoBack" bookmark.
$endRng = $d.Range($d.Content.End, $d.Content.End)
$newParasXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:b/><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:u w:val="single"/></w:rPr><w:t xml:space="preserve">Central Limit Theorem (CLT) and t-distribution </w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">If we apply the CLT, the distribution of the t-statistic is normal with mean 0 and standard deviation 1. </w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>
'@
$endRng.InsertXML($newParasXml)
